$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The elevator repair items will now always spawn on the floor, so mark the
# "repair items on location" task, the related lift-breakage task and the
# enemy-scaling-per-floor task as done (row 1, row 3, row 16).
# Copy the formatting (fill/border/font) already used by the other "Done"
# cells (e.g. B2) so these cells pick up the same style instead of a new one.
$ws.Range("B2").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("B1").Value = "Cделано"
$ws.Range("B3").Value = "Cделано"
$ws.Range("B16").Value = "Cделано"

# Move the active selection to B3 and scroll the sheet so row 5 is at the
# top-left of the view (topLeftCell = A5).
$ws.Range("B3").Select()
$win = $ws.Application.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 5
